$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43: fill in the previously-empty "tijd" cell for the Week 15 header row ---
$ws.Range("B43").Value = "7 uur 30 minuten"

# --- Row 46: complete the existing row (only Q46 had data before) ---
$ws.Range("A46").Value = 43611
$ws.Range("A46").NumberFormat = "d-mmm"
$ws.Range("B46").Value = "1 uur 20 minuten"
$ws.Range("C46").Value = "elapsed tijd locatie debuggen"

# --- Row 47: new "Week 16" header row ---
$ws.Range("A47").Value = "Week 16"
$ws.Range("A47").Font.Bold = $true

$ws.Hyperlinks.Add($ws.Range("Q47"), "https://stackoverflow.com/questions/4803248/dim-screen-while-user-inactive")
$ws.Range("Q47").Value = "https://stackoverflow.com/questions/4803248/dim-screen-while-user-inactive"
$ws.Range("Q47").Style = "Hyperlink"

# --- Row 48: new data row ---
$ws.Range("A48").Value = 43612
$ws.Range("A48").NumberFormat = "d-mmm"
$ws.Range("B48").Value = "30 minuten"
$ws.Range("C48").Value = "Project structuur documentatie maken"

$ws.Hyperlinks.Add($ws.Range("Q48"), "https://developer.android.com/training/graphics/opengl/touch", "java", "", "https://developer.android.com/training/graphics/opengl/touch - java")
$ws.Range("Q48").Value = "https://developer.android.com/training/graphics/opengl/touch#java"
$ws.Range("Q48").Style = "Hyperlink"

# --- Row 49: new data row (no Q entry) ---
$ws.Range("A49").Value = 43614
$ws.Range("A49").NumberFormat = "d-mmm"
$ws.Range("B49").Value = "2 uur 15 minuten"
$ws.Range("C49").Value = "Documentatie afwerken en nalezen, informatie scherm donker maken opzoeken en toepassen in project en debuggen"

# --- Update selection / scroll position to match the author's final view ---
$ws.Range("A37").Select() | Out-Null
$ws.Range("C49").Select() | Out-Null
